# Adding new testcases for iAuthor and changing smoke suite of iProctor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "users" to "Worksheet"
$ws.Name = "Worksheet"

# Refresh the generated test-data values in row 2 (new random creds/ids)
$ws.Range("A2").Value = "ZhyoF284"
$ws.Range("B2").Value = 23100602
$ws.Range("C2").Value = "qfzzgjw49"
$ws.Range("D2").Value = "A&3sP!h6"
$ws.Range("F2").Value = "xgWJZTBa"
$ws.Range("G2").Value = "ydGc"

# Strip the bold header / bordered cell formatting back to plain default styling
[void]$ws.Range("A1:H2").ClearFormats()

# Collapse the selection down to the single top-left cell
[void]$ws.Range("A1").Select()
